# Update bracket pick cells on the "Bracket" sheet to reflect the
# winners advancing through the tournament rounds.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$updates = @{
    "D2"  = "Okapi"
    "E4"  = "Okapi"
    "M4"  = "Puffer Fish"
    "F8"  = "Okapi"
    "N10" = "Trapdoor Spider"
    "E12" = "Striped dolphin"
    "M12" = "Lungfish"
    "D14" = "Striped dolphin"
    "N14" = "Lungfish"
    "G16" = "Badger Bat"
    "F24" = "Badger Bat"
    "H32" = "Rock Hyrax "
    "M36" = "Emperor Penguin"
    "N42" = "Peacock Wrasse"
    "E44" = "Mara"
    "D46" = "Mara"
    "G48" = "Rock Hyrax "
    "D50" = "Itjaritjari"
    "N50" = "Spotted sandpiper"
    "E52" = "Dik Dik"
    "M52" = "Giant Water Bug"
    "D54" = "Dik Dik"
    "N54" = "Giant Water Bug"
    "F56" = "Rock Hyrax "
    "E60" = "Rock Hyrax "
    "D62" = "Rock Hyrax "
    "N62" = "Three-Spined stickleback "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
